$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I20").Value = -0.9712759750021939
$ws.Range("J20").Value = 0.2859287179783679
$ws.Range("K20").Value = 0.5190875794184199
$ws.Range("L20").Value = 2.848823562284922
